# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.897.32"
$ws.Range("E2").Value = "'  +2.77%  "
$ws.Range("D3").Value = "'2.453.60"
$ws.Range("E3").Value = "'  +2.02%  "
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'575.34"
$ws.Range("E5").Value = "'  +1.29%  "
$ws.Range("D6").Value = "'146.15"
$ws.Range("E6").Value = "'  +2.60%  "
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("E8").Value = "'  +0.91%  "
$ws.Range("D9").Value = "'2.453.13"
$ws.Range("E9").Value = "'  +1.59%  "
$ws.Range("E10").Value = "'  +2.26%  "
$ws.Range("E11").Value = "'  +2.48%  "
$ws.Range("D12").Value = "'5.27"
$ws.Range("E12").Value = "'  +1.08%  "
$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "'  +2.34%  "
$ws.Range("D14").Value = "'28.58"
$ws.Range("E14").Value = "'  +7.98%  "
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("E15").Value = "'  +3.92%  "
$ws.Range("D16").Value = "'2.898.67"
$ws.Range("E16").Value = "'  +3.64%  "
$ws.Range("D17").Value = "'62.992.51"
$ws.Range("E17").Value = "'  +3.36%  "
$ws.Range("D18").Value = "'2.456.55"
$ws.Range("E18").Value = "'  +2.00%  "
$ws.Range("D19").Value = "'7.92"
$ws.Range("E19").Value = "'  -1.32%  "
$ws.Range("D20").Value = "'11.05"
$ws.Range("E20").Value = "'  +3.07%  "
$ws.Range("D21").Value = "'329.92"
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "'  +1.09%  "
$ws.Range("E23").Value = "'  +8.88%  "
$ws.Range("E24").Value = "'  +0.01%  "
$ws.Range("D25").Value = "'66.44"
$ws.Range("E25").Value = "'  +1.99%  "
$ws.Range("B26").Value = "'Bittensor"
$ws.Range("C26").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D26").Value = "'654.50"
$ws.Range("E26").Value = "'  +7.76%  "
$ws.Range("B27").Value = "'Binance-PegBSC-USD"
$ws.Range("C27").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.21"
$ws.Range("E27").Value = "'  +20.63%  "
$ws.Range("D28").Value = "'8.59"
$ws.Range("E28").Value = "'  +3.29%  "
$ws.Range("D29").Value = "'0.0₃0992"
$ws.Range("E29").Value = "'  +4.65%  "
$ws.Range("E30").Value = "'  +2.56%  "
$ws.Range("D31").Value = "'8.20"
$ws.Range("E31").Value = "'  +2.54%  "
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "'  +3.90%  "
$ws.Range("E33").Value = "'  +3.82%  "
$ws.Range("D34").Value = "'0.138"
$ws.Range("E34").Value = "'  +4.86%  "
$ws.Range("B35").Value = "'BabyDogeCoin"
$ws.Range("C35").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D35").Value = "'0.0₆0377"
$ws.Range("E35").Value = "'  +33.50%  "
$ws.Range("B36").Value = "'ImmutableX"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.48"
$ws.Range("E36").Value = "'  +0.53%  "
$ws.Range("B37").Value = "'FirstDigitalUSD"
$ws.Range("C37").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  +0.14%  "
$ws.Range("B38").Value = "'NEARProtocol"
$ws.Range("C38").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.78"
$ws.Range("E38").Value = "'  +3.36%  "
$ws.Range("B39").Value = "'RenderToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'5.50"
$ws.Range("E39").Value = "'  +3.75%  "
$ws.Range("B40").Value = "'PolygonEcosystemToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.373"
$ws.Range("E40").Value = "'  -0.11%  "
$ws.Range("D41").Value = "'152.76"
$ws.Range("E41").Value = "'  +1.10%  "
$ws.Range("B42").Value = "'EthereumClassic"
$ws.Range("C42").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D42").Value = "'18.80"
$ws.Range("E42").Value = "'  +2.41%  "
$ws.Range("B43").Value = "'dogwifhat"
$ws.Range("C43").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.70"
$ws.Range("E43").Value = "'  +7.00%  "
$ws.Range("B44").Value = "'Stacks"
$ws.Range("C44").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.76"
$ws.Range("E44").Value = "'  +3.53%  "
$ws.Range("B45").Value = "'OKB"
$ws.Range("C45").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'42.54"
$ws.Range("E45").Value = "'  +1.45%  "
$ws.Range("B46").Value = "'USDe"
$ws.Range("C46").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("B47").Value = "'WhiteBITCoin"
$ws.Range("C47").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'14.99"
$ws.Range("E47").Value = "'  +27.10%  "
$ws.Range("B48").Value = "'Aave"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'146.36"
$ws.Range("E48").Value = "'  +3.62%  "
$ws.Range("B49").Value = "'Filecoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'3.62"
$ws.Range("E49").Value = "'  +2.62%  "
$ws.Range("B50").Value = "'InjectiveProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'20.59"
$ws.Range("E50").Value = "'  +3.45%  "
$ws.Range("B51").Value = "'Mantle"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.606"
$ws.Range("E51").Value = "'  +2.23%  "
